$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 160, shifting existing rows 160..267 down to 161..268
$ws.Rows("160:160").Insert()

# Populate the new row 160 with the same values as the (now shifted) row 161,
# except for the Fecha (D) and Volumen (M) columns which get new values.
$ws.Range("A160").Value = 10
$ws.Range("B160").Value = "Vega Modelo de Temuco"
$ws.Range("C160").Value = "La Araucanía"
$ws.Range("D160").Value = 45126
$ws.Range("E160").Value = 9
$ws.Range("F160").Value = "Fruta"
$ws.Range("G160").Value = 100104
$ws.Range("H160").Value = "Frutos de pepita"
$ws.Range("I160").Value = 100104001
$ws.Range("J160").Value = "Granada"
$ws.Range("K160").Value = "Wonderfull"
$ws.Range("L160").Value = "Primera"
$ws.Range("M160").Value = 95
$ws.Range("N160").Value = 14000
$ws.Range("O160").Value = 14000
$ws.Range("P160").Value = 14000
$ws.Range("Q160").Value = "$/bandeja 10 kilos granel"
$ws.Range("R160").Value = "Provincia de Limarí"
$ws.Range("S160").Value = 1400
$ws.Range("T160").Value = 10

$ws.Range("D160").NumberFormat = "YYYY-MM-DD HH:MM:SS"
